$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = 1.334851980209351
$ws.Range("B1").Value = 2.748502969741821
$ws.Range("C1").Value = 2.905020475387573
$ws.Range("D1").Value = 1.480071902275085
$ws.Range("E1").Value = 1.079721450805664
